# Insert a new weekly data record at row 243, shifting existing rows
# 243-327 down to 244-328 (dimension grows from A1:R327 to A1:R328).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(243).Insert()

$ws.Cells.Item(243,1).Value()  = 7
$ws.Cells.Item(243,2).Value()  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(243,3).Value()  = "Ñuble"
$ws.Cells.Item(243,4).Value()  = 44524
$ws.Cells.Item(243,5).Value()  = 16
$ws.Cells.Item(243,6).Value()  = 100112020
$ws.Cells.Item(243,7).Value()  = "Tomate"
$ws.Cells.Item(243,8).Value()  = "Larga vida"
$ws.Cells.Item(243,9).Value()  = "Primera"
$ws.Cells.Item(243,10).Value() = 400
$ws.Cells.Item(243,11).Value() = 7000
$ws.Cells.Item(243,12).Value() = 8000
$ws.Cells.Item(243,13).Value() = 7500
$ws.Cells.Item(243,14).Value() = "`$/caja 15 kilos"
$ws.Cells.Item(243,15).Value() = "Región del Maule"
$ws.Cells.Item(243,16).Value() = 500
$ws.Cells.Item(243,17).Value() = 15
$ws.Cells.Item(243,18).Value() = "Hortaliza"

$ws.Cells.Item(243,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
